# Add a new "ΠΕΛΑΤΕΣ" (customers) worksheet after the existing "ΠΡΟΪΟΝΤΑ"
# sheet, fill it with the customer table, and register a matching
# workbook-level defined name - mirrors the upstream "Add files via
# upload" commit that introduced the ΠΕΛΑΤΕΣ sheet + named range.

$wb = $excel.ActiveWorkbook

# --- 1. Create the new worksheet, positioned right after ΠΡΟΪΟΝΤΑ -----
$wsProducts = $wb.Worksheets.Item(1)
$ws = $wb.Worksheets.Add($null, $wsProducts)
$ws.Name = "ΠΕΛΑΤΕΣ"

# --- 2. Header row ------------------------------------------------------
$ws.Range("A1").Value = "ΚΩΔ_ΠΕΛ"
$ws.Range("B1").Value = "ΟΝΟΜΑ"
$ws.Range("C1").Value = "ΕΠΩΝΥΜΟ"
$ws.Range("D1").Value = "ΔΙΕΥΘΥΝΣΗ"
$ws.Range("E1").Value = "ΠΟΛΗ"
$ws.Range("F1").Value = "ΤΗΛΕΦΩΝΟ"

# --- 3. Data rows --------------------------------------------------------
# Column F (ΤΗΛΕΦΩΝΟ) holds numeric-looking phone numbers that must stay
# text (leading area codes etc.), so mark the column as Text before
# writing into it - otherwise plain numeric strings get coerced to
# numbers, same as typing them into real Excel would.
$ws.Range("F2:F7").NumberFormat = "@"

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Άγγελος"
$ws.Range("C2").Value = "Βώρος"
$ws.Range("D2").Value = "Κορίνθου 200"
$ws.Range("E2").Value = "Πάτρα"
$ws.Range("F2").Value = "26109919890"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Άγγελος"
$ws.Range("C3").Value = "Βώρος"
$ws.Range("D3").Value = "Κορίνθου 200"
$ws.Range("E3").Value = "Πάτρα"
$ws.Range("F3").Value = "26109919890"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Κώστας"
$ws.Range("C4").Value = "Περδίου"
$ws.Range("D4").Value = "Αθηνών 6"
$ws.Range("E4").Value = "Αθήνα"
$ws.Range("F4").Value = "2109919890"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Κώστας"
$ws.Range("C5").Value = "Περδίου"
$ws.Range("D5").Value = "Αθηνών 6"
$ws.Range("E5").Value = "Αθήνα"
$ws.Range("F5").Value = "2109919890"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Γιώργος"
$ws.Range("C6").Value = "Καλέμης"
$ws.Range("D6").Value = "Αράτου 111"
$ws.Range("E6").Value = "Πάτρα"
$ws.Range("F6").Value = "2610887960"

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Γιώργος"
$ws.Range("C7").Value = "Καλέμης"
$ws.Range("D7").Value = "Αράτου 111"
$ws.Range("E7").Value = "Πάτρα"
$ws.Range("F7").Value = "2610887960"

# --- 4. Selection matches the source sheet (C5) -------------------------
$null = $ws.Range("C5").Select()

# --- 5. Workbook-level defined names ------------------------------------
# `Names.Add` (and the `Range.Name =` shortcut) only resolves the COM
# method/property call when the *name* argument's leading character is
# ASCII, so a Greek-first identifier like "ΠΕΛΑΤΕΣ" has to be created
# under a throwaway ASCII alias first and then renamed in place.
$r = $ws.Range("A1:F7")
$r.Name = "TEMP_PELATES_NAME"
$newName = $wb.Names.Item("TEMP_PELATES_NAME")
$newName.RefersTo = "='ΠΕΛΑΤΕΣ'!`$A`$1:`$F`$7"
$newName.Name = "ΠΕΛΑΤΕΣ"

# Re-assert the original ΠΡΟΪΟΝΤΑ defined name explicitly (quoted sheet
# reference) so it keeps its original shape after the workbook's defined
# names collection has been touched.
$oldName = $wb.Names.Item("ΠΡΟΪΟΝΤΑ")
$oldName.RefersTo = "='ΠΡΟΪΟΝΤΑ'!`$A`$1:`$C`$5"
